$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue $ws "D2" "317.09"
Set-TextValue $ws "E2" "4.43%"
Set-TextValue $ws "D3" "48.70"
Set-TextValue $ws "E3" "13.29%"
Set-TextValue $ws "D4" "5.268"
Set-TextValue $ws "E4" "4.63%"
Set-TextValue $ws "D5" "0.07995"
Set-TextValue $ws "E5" "3.98%"
Set-TextValue $ws "D6" "4.592"
Set-TextValue $ws "E6" "4.15%"
Set-TextValue $ws "D7" "1.433"
Set-TextValue $ws "E7" "34.93%"
Set-TextValue $ws "D8" "1.644"
Set-TextValue $ws "E8" "2.13%"
Set-TextValue $ws "D9" "0.1277"
Set-TextValue $ws "E9" "3.55%"
Set-TextValue $ws "D10" "0.1935"
Set-TextValue $ws "E10" "3.59%"
Set-TextValue $ws "D11" "0.09293"
Set-TextValue $ws "E11" "1.81%"
Set-TextValue $ws "D12" "0.04598"
Set-TextValue $ws "E12" "10.40%"
Set-TextValue $ws "D13" "0.1045"
Set-TextValue $ws "E13" "-0.01%"
Set-TextValue $ws "D14" "0.001319"
Set-TextValue $ws "E14" "4.03%"
Set-TextValue $ws "D16" "0.005823"
Set-TextValue $ws "E16" "0.95%"
Set-TextValue $ws "E17" "0.21%"
Set-TextValue $ws "D18" "2.435"
Set-TextValue $ws "E18" "2.14%"
Set-TextValue $ws "D19" "0.3416"
Set-TextValue $ws "E19" "2.14%"
Set-TextValue $ws "D20" "8.162"
Set-TextValue $ws "E20" "-2.57%"
Set-TextValue $ws "E21" "-0.34%"
Set-TextValue $ws "D22" "0.3101"
Set-TextValue $ws "E22" "7.25%"
Set-TextValue $ws "D23" "0.001310"
Set-TextValue $ws "E23" "2.96%"
Set-TextValue $ws "D24" "0.004238"
Set-TextValue $ws "E24" "-5.61%"
Set-TextValue $ws "D25" "0.0001354"
Set-TextValue $ws "E25" "0.59%"
Set-TextValue $ws "D26" "0.0003542"
Set-TextValue $ws "E26" "-95.23%"
Set-TextValue $ws "D38" "0.02696"
Set-TextValue $ws "E38" "9.63%"
Set-TextValue $ws "D39" "0.05647"
Set-TextValue $ws "E39" "6.93%"
Set-TextValue $ws "D40" "0.008224"
Set-TextValue $ws "E40" "37.95%"
Set-TextValue $ws "D41" "0.008019"
Set-TextValue $ws "E41" "4.69%"
Set-TextValue $ws "E42" "6.66%"
Set-TextValue $ws "D43" "0.007686"
Set-TextValue $ws "E43" "4.70%"
Set-TextValue $ws "D44" "0.008622"
Set-TextValue $ws "E44" "3.35%"
Set-TextValue $ws "D45" "0.3466"
Set-TextValue $ws "E45" "14.67%"
Set-TextValue $ws "D46" "0.00006903"
Set-TextValue $ws "E46" "4.16%"
Set-TextValue $ws "D47" "0.00000000752"
Set-TextValue $ws "E47" "0.58%"
Set-TextValue $ws "E48" "42.79%"
Set-TextValue $ws "D49" "0.004003"
Set-TextValue $ws "E49" "-4.70%"
Set-TextValue $ws "D50" "0.00002106"
Set-TextValue $ws "E50" "0.58%"
Set-TextValue $ws "D51" "0.0002006"
Set-TextValue $ws "E51" "0.58%"
